$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E header: "2025-03-31" must stay literal text, not get
# auto-converted to a date serial value. Enter it as a formula-text literal,
# then convert the cell to its computed value in place (values-only paste),
# which keeps the General style (no NumberFormat mutation => no new style).
$ws.Range("E1").Formula = '="2025-03-31"'
$ws.Range("E1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# New student rows
$ws.Range("A4").Value = "Aaron"
$ws.Range("B4").Value = "Javier"
$ws.Range("C4").Value = 3

$ws.Range("A5").Value = "Gael"
$ws.Range("B5").Value = "Franco"
$ws.Range("C5").Value = 4

$ws.Range("A6").Value = "Prueba"
$ws.Range("B6").Value = "uno"
$ws.Range("C6").Value = 5

# Copy the green attendance-cell formatting from the existing D column onto
# the new attendance cells, so the existing style is reused (not re-created).
$ws.Range("D3").Copy() | Out-Null
$ws.Range("E2:E4").PasteSpecial(-4122) | Out-Null
$ws.Range("D4:D5").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
